# Sync the Listeria table with its Markdown source:
#   - Row 20 (1999-2000, Osaka) and Row 21 (1998-2001, Nissui) are removed.
#   - The former Row 22 (1988-2004, NIHS) becomes the new Row 20.
# Net effect: old row 22's values move up into row 20, then rows 21:22 are deleted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the last row's (row 22) values before shifting anything.
$a22 = $ws.Range("A22").Value()
$b22 = $ws.Range("B22").Value()
$c22 = $ws.Range("C22").Value()

# Overwrite row 20 with what used to be row 22.
$ws.Range("A20").Value = $a22
$ws.Range("B20").Value = $b22
$ws.Range("C20").Value = $c22

# Remove the now-duplicated old rows 21 and 22, shrinking the table to A1:C20.
$ws.Rows("21:22").Delete()
